# Commit: "update ip to 192.168.1.x"
# Renumber the lab's IP addresses from the 192.168.121.0/24 range to
# 192.168.1.0/24, and leave the selection where the author left it when
# they saved the file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column F holds the "IP" values (row 1 is the header).
$ws.Range("F2").Value = "192.168.1.99"    # catpc           was 192.168.121.100
$ws.Range("F3").Value = "192.168.1.111"   # k8smaster       was 192.168.121.111
$ws.Range("F4").Value = "192.168.1.112"   # k8sworker1      was 192.168.121.112
$ws.Range("F5").Value = "192.168.1.113"   # k8sworker2      was 192.168.121.113
$ws.Range("F6").Value = "192.168.1.114"   # rancherserver   was 192.168.121.114
$ws.Range("F9").Value = "192.168.1.100"   # loadbalance     was 192.168.121.99

# Restore the author's last selection/active cell on Sheet1.
$ws.Activate()
$ws.Range("F13").Select()

# Match the saved window size from the workbook view (best effort - the
# host window chrome size isn't always persisted by every engine).
$excel.ActiveWindow.Width = 21000
$excel.ActiveWindow.Height = 7902
